# Fix complicated issue with death zones
# See "Bug 2" text under Part B Notes.txt for more details.
#
# On the Board sheet, the last mini-board (columns AT:BA, rows 2:9)
# had a bunch of stray piece markers left over from testing. Clear
# them out and correct a couple of cells so the board reflects the
# real (fixed) state, then highlight the piece that sits in the
# "death zone" in red so it's easy to spot.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Board")

# Cells that should go back to being blank (but keep their existing
# border/fill formatting - i.e. just clear the contents).
$cellsToClear = @(
    "AU2", "AV2", "AY2",
    "AU3", "AV3", "BA3",
    "AT4", "AU4", "AW4",
    "AV5", "AX5", "BA5",
    "AT6", "AV6", "BA6",
    "AT7", "AX7", "AY7",
    "AV8", "AZ8",
    "AW9"
)
foreach ($addr in $cellsToClear) {
    $ws.Range($addr).ClearContents()
}

# BB2 was a stray numeric helper value outside the real board - remove
# the cell completely (shrinks the sheet's used range back to BA26).
$ws.Range("BB2").Clear()

# A few cells were showing the wrong piece / were missing their piece.
$ws.Range("AZ3").Value = "X"
$ws.Range("AZ4").Value = "B"
$ws.Range("AZ5").Value = "W"

# AY5 holds the piece that is actually in danger in the death zone -
# call it out with a red font (new style, same fill/border as the
# other "Input" board cells).
$ws.Range("AY5").Value = "W"
$ws.Range("AY5").Font.Color = 255

# Reflect where the user was last working when they made this fix.
$ws.Range("AW2").Select()
